$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 15:03"

# Row 5
$ws.Range("B5").Value = 250561
$ws.Range("C5").Value = 2260
$ws.Range("D5").Value = 154718
$ws.Range("E5").Value = 70230
$ws.Range("F5").Value = 2254
$ws.Range("G5").Value = 185
$ws.Range("H5").Value = 25613

# Row 19
$ws.Range("B19").Value = 41087
$ws.Range("C19").Value = 317
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 35669
$ws.Range("F19").Value = 683
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = 5168

# Row 21 (country -> Arabia Saudita)
$ws.Range("A21").Value = "Arabia Saudita"
$ws.Range("B21").Value = 30251
$ws.Range("C21").Value = 1595
$ws.Range("D21").Value = 5431
$ws.Range("E21").Value = 24620
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 200

# Row 22 (country -> Suiza)
$ws.Range("A22").Value = "Suiza"
$ws.Range("B22").Value = 30009
$ws.Range("C22").Value = 28
$ws.Range("D22").Value = 25200
$ws.Range("E22").Value = 3019
$ws.Range("F22").Value = 141
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 1790

# Row 23
$ws.Range("B23").Value = 25702
$ws.Range("C23").Value = 178
$ws.Range("D23").Value = 1743
$ws.Range("E23").Value = 22885
$ws.Range("F23").Value = 134
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 1074

# Row 42
$ws.Range("D42").Value = 7296
$ws.Range("E42").Value = 2022
$ws.Range("F42").Value = 49
$ws.Range("G42").Value = 10
$ws.Range("H42").Value = 503

# Row 55
$ws.Range("E55").Value = 1666
$ws.Range("F55").Value = 48
$ws.Range("G55").Value = 6
$ws.Range("H55").Value = 246

# Row 72
$ws.Range("B72").Value = 2204
$ws.Range("C72").Value = 15
$ws.Range("D72").Value = 1454
$ws.Range("E72").Value = 740
$ws.Range("F72").Value = 8

# Row 73 (country -> Croacia)
$ws.Range("A73").Value = "Croacia"
$ws.Range("B73").Value = 2112
$ws.Range("C73").Value = 11
$ws.Range("D73").Value = 1560
$ws.Range("E73").Value = 469
$ws.Range("F73").Value = 14
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 83

# Row 74 (country -> Camerun)
$ws.Range("A74").Value = "Camerun"
$ws.Range("B74").Value = 2104
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 953
$ws.Range("E74").Value = 1087
$ws.Range("F74").Value = 12
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 64

# Row 83
$ws.Range("B83").Value = 1526
$ws.Range("C83").Value = 8
$ws.Range("D83").Value = 1013
$ws.Range("E83").Value = 427
$ws.Range("F83").Value = 21
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 86

# Row 98 (country -> Sudan)
$ws.Range("A98").Value = "Sudan"
$ws.Range("B98").Value = 778
$ws.Range("C98").Value = 100
$ws.Range("D98").Value = 70
$ws.Range("E98").Value = 663
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = 45

# Row 99 (country -> Somalia)
$ws.Range("A99").Value = "Somalia"
$ws.Range("B99").Value = 756
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 61
$ws.Range("E99").Value = 660
$ws.Range("F99").Value = 2
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 35

# Row 100 (country -> Sri Lanka)
$ws.Range("A100").Value = "Sri Lanka"
$ws.Range("C100").Value = 4
$ws.Range("D100").Value = 197
$ws.Range("E100").Value = 549
$ws.Range("F100").Value = 1
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 9

# Row 101 (country -> Niger)
$ws.Range("A101").Value = "Niger"
$ws.Range("B101").Value = 755
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 534
$ws.Range("E101").Value = 184
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 37

# Row 102 (country -> Principado de Andorra)
$ws.Range("A102").Value = "Principado de Andorra"
$ws.Range("B102").Value = 750
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 499
$ws.Range("E102").Value = 206
$ws.Range("F102").Value = 16
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 45

# Row 103 (country -> Costa Rica)
$ws.Range("A103").Value = "Costa Rica"
$ws.Range("B103").Value = 742
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 399
$ws.Range("E103").Value = 337
$ws.Range("F103").Value = 5
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 6

# Row 104 (country -> Libano)
$ws.Range("A104").Value = "Libano"
$ws.Range("B104").Value = 741
$ws.Range("C104").Value = 1
$ws.Range("D104").Value = 201
$ws.Range("E104").Value = 515
$ws.Range("F104").Value = 43
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 25

# Row 105 (country -> Guatemala)
$ws.Range("A105").Value = "Guatemala"
$ws.Range("B105").Value = 730
$ws.Range("C105").Value = 27
$ws.Range("D105").Value = 79
$ws.Range("E105").Value = 632
$ws.Range("F105").Value = 5
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 19

# Row 106 (country -> Crucero)
$ws.Range("A106").Value = "Crucero"
$ws.Range("B106").Value = 712
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 645
$ws.Range("E106").Value = 54
$ws.Range("F106").Value = 4
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 13

# Row 107 (country -> Mayotte)
$ws.Range("A107").Value = "Mayotte"
$ws.Range("B107").Value = 686
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 352
$ws.Range("E107").Value = 328
$ws.Range("F107").Value = 6
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 6

# Row 108 (country -> Consejo Danes para los Refugiados)
$ws.Range("A108").Value = "Consejo Danes para los Refugiados"
$ws.Range("B108").Value = 682
$ws.Range("C108").Value = 0
$ws.Range("D108").Value = 80
$ws.Range("E108").Value = 568
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 34

# Row 116
$ws.Range("B116").Value = 535
$ws.Range("C116").Value = 45
$ws.Range("D116").Value = 182
$ws.Range("E116").Value = 329
$ws.Range("F116").Value = 2

# Row 199 (country -> Burundi)
$ws.Range("A199").Value = "Burundi"
$ws.Range("D199").Value = 7
$ws.Range("E199").Value = 7
$ws.Range("H199").Value = 1

# Row 200 (country -> San Cristobal y Nieves)
$ws.Range("A200").Value = "San Cristobal y Nieves"
$ws.Range("D200").Value = 8
$ws.Range("E200").Value = 7
$ws.Range("H200").Value = 0
